$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Slide 16 table: switch the graphicFrame's table style to the new style
#    ({3286A7B5-...} -> {A5FD4C48-...}), matching the tableStyleId change in
#    the target XML diff.
# ---------------------------------------------------------------------------
$slide16 = $p.Slides.Item(16)
for ($i = 1; $i -le $slide16.Shapes.Count; $i++) {
    $shp = $slide16.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{A5FD4C48-E46B-46F7-960E-247DF1D8B847}", $true)
    }
}

# ---------------------------------------------------------------------------
# 2) Theme swap: the deck's live theme (currently the "Integral" palette)
#    gets recoloured to the stock "Office Theme" palette. All slides share
#    the one slide master/theme, so touching any slide's ThemeColorScheme
#    updates that single shared theme part.
#    Order: dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink (indices 1-12).
# ---------------------------------------------------------------------------
$officeColors = @(0, 16777215, 6968388, 15132391, 13998939, 3243501, 10855845, 49407, 12874308, 4697456, 12673797, 7491477)

$slide1 = $p.Slides.Item(1)
$tcs = $slide1.ThemeColorScheme
for ($i = 1; $i -le $tcs.Count; $i++) {
    $tcs.Colors($i).RGB = $officeColors[$i - 1]
}
